$d = $word.ActiveDocument
$bullet = [char]8226

function Get-ParaIndexByText([string]$text) {
    # Locate a paragraph by its (unique) text and return its 1-based index
    # in $d.Paragraphs, using Find rather than a hard-coded number.
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r.Paragraphs.Item(1).Index
}

function Add-BulletAfter([int]$index, [string]$text) {
    # Insert a new plain ("Normal") paragraph right after $d.Paragraphs.Item($index)
    # containing $text, and return the new paragraph's index.
    $after = $d.Paragraphs.Item($index)
    $after.Range.InsertParagraphAfter() | Out-Null
    $new = $d.Paragraphs.Item($index + 1)
    if ($new.Style.NameLocal -ne "Normal") {
        $new.Style = "Normal"
    }
    $new.Range.Text = $text
    return $index + 1
}

# --- 1. Update "Software Development and Innovation" heading text ---
$d.Content.Find.Execute(
    "Software Development and Innovation", $true, $false, $false, $false, $false,
    $true, 1, $false, "Technical Innovation & Platform Development", 2) | Out-Null

# --- 2. Expand the first bullet under that heading ---
$d.Content.Find.Execute(
    "$bullet Conceived and deployed redistricting software used by thousands of analysts nationwide",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "$bullet Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide",
    2) | Out-Null

# --- 3. Insert four new bullets after it ---
$idx = Get-ParaIndexByText("$bullet Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide")
$idx = Add-BulletAfter $idx "$bullet Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
$idx = Add-BulletAfter $idx "$bullet Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors"
$idx = Add-BulletAfter $idx "$bullet Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls"
$idx = Add-BulletAfter $idx "$bullet Created SimCrisis platform for humanitarian intervention modeling used by International Red Cross and UNICEF"

# --- 4. "Developed boundary estimation ..." is left untouched ---

# --- 5. Turn "Created econometric simulation platform ..." into the new Heading3 ---
$idx = Get-ParaIndexByText("$bullet Created econometric simulation platform for humanitarian intervention modeling")
$p = $d.Paragraphs.Item($idx)
$p.Style = "Heading 3"
$p.Range.Text = "Data Engineering & Analytics"

# --- 6. Insert six new bullets after the new heading ---
$idx = Add-BulletAfter $idx "$bullet Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change"
$idx = Add-BulletAfter $idx "$bullet Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%"
$idx = Add-BulletAfter $idx "$bullet Developed advanced data pipelines for machine learning applications enhancing consumer segmentation and predictive modeling"
$idx = Add-BulletAfter $idx "$bullet Built fraud detection systems for campaign finance data analysis across multi-terabyte datasets"
$idx = Add-BulletAfter $idx "$bullet Transformed small data team into big data engineering team using Hadoop Clusters and Hive on AWS"
$idx = Add-BulletAfter $idx "$bullet Introduced version control and Agile methodologies, improving project delivery timelines by 40%"

# --- 7. Insert the "Research Leadership & Client Success" Heading3 after the last bullet ---
$after = $d.Paragraphs.Item($idx)
$after.Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Style = "Heading 3"
$p.Range.Text = "Research Leadership & Client Success"

# --- 8. Insert three new bullets after it ---
$idx = Add-BulletAfter $idx "$bullet Led multi-million dollar research projects involving sensitive consumer data with privacy compliance"
$idx = Add-BulletAfter $idx "$bullet Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders"
$idx = Add-BulletAfter $idx "$bullet Delivered actionable consumer insights and market intelligence for political candidates and major organizations"

# --- 9. "Built comprehensive survey operations ..." is left untouched ---

# --- 10. Append two final bullets after it ---
$idx = Get-ParaIndexByText("$bullet Built comprehensive survey operations platform from RFP through deployment")
$idx = Add-BulletAfter $idx "$bullet Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership"
$idx = Add-BulletAfter $idx "$bullet Redistricting analysis used in court cases with rigorous methodology and expert testimony"

Write-Output "done"
